# Update the "Förändrad" (Changed) date column (column C) for every data
# row on the sheet: the stored serial date value 45177 becomes 45178.
# Data occupies rows 2-471 (row 1 is the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C471").Value = 45178
